$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 516.1818
$ws.Range("J19").Value = 430.72726
$ws.Range("L19").Value = 430.72726
$ws.Range("N19").Value = -780.72726
$ws.Range("H40").Value = 7437.9165
$ws.Range("I40").Value = 1299.5
$ws.Range("J40").Value = 8665.6
$ws.Range("K40").Value = 1299.5
$ws.Range("L40").Value = 8665.6
$ws.Range("M40").Value = -1124.5
$ws.Range("N40").Value = -9015.6
$ws.Range("H112").Value = 2071.0417
$ws.Range("J112").Value = 2100.2273
$ws.Range("L112").Value = 6300.6819
$ws.Range("N112").Value = -8516.6819

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H61").Value = 6701.8823
$ws.Range("I61").Value = 5770.4443
$ws.Range("K61").Value = 5770.4443
$ws.Range("M61").Value = -5558.4443
$ws.Range("H122").Value = 4678.7646
$ws.Range("I122").Value = 2570.724
$ws.Range("J122").Value = 16905.4
$ws.Range("K122").Value = 7712.172
$ws.Range("L122").Value = 50716.2
$ws.Range("M122").Value = -5262.172
$ws.Range("N122").Value = -55616.2
$ws.Range("H132").Value = 7071.6875
$ws.Range("I132").Value = 7454.8
$ws.Range("J132").Value = 6433.1665
$ws.Range("K132").Value = 22364.4
$ws.Range("L132").Value = 19299.4995
$ws.Range("M132").Value = -19834.4
$ws.Range("N132").Value = -24359.4995
$ws.Range("H136").Value = 6701.8823
$ws.Range("I136").Value = 5770.4443
$ws.Range("K136").Value = 17311.3329
$ws.Range("M136").Value = -14761.3329

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5668.4375
$ws.Range("I134").Value = 3997.889
$ws.Range("K134").Value = 11993.667
$ws.Range("M134").Value = -9458.667000000001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2985.39
$ws.Range("J31").Value = 3920.08
$ws.Range("L31").Value = 3920.08
$ws.Range("N31").Value = -4510.08
$ws.Range("H34").Value = 2985.39
$ws.Range("J34").Value = 3920.08
$ws.Range("L34").Value = 3920.08
$ws.Range("N34").Value = -4324.08
$ws.Range("H99").Value = 10692951
$ws.Range("J99").Value = 15393439
$ws.Range("L99").Value = 15393439
$ws.Range("N99").Value = -15396435
$ws.Range("H126").Value = 10692951
$ws.Range("J126").Value = 15393439
$ws.Range("L126").Value = 46180317
$ws.Range("N126").Value = -46185257

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 223369
$ws.Range("J37").Value = 223369
$ws.Range("L37").Value = 670107
$ws.Range("N37").Value = -670331
$ws.Range("H50").Value = 1256.125
$ws.Range("J50").Value = 1919.8
$ws.Range("L50").Value = 5759.4
$ws.Range("N50").Value = -6721.4
$ws.Range("H53").Value = 1256.125
$ws.Range("J53").Value = 1919.8
$ws.Range("L53").Value = 5759.4
$ws.Range("N53").Value = -6721.4
$ws.Range("H114").Value = 1450.9231
$ws.Range("I114").Value = 473.5
$ws.Range("K114").Value = 1420.5
$ws.Range("M114").Value = 1833.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 50000732
$ws.Range("I97").Value = 636.3333
$ws.Range("K97").Value = 636.3333
$ws.Range("M97").Value = -140.3333
$ws.Range("H102").Value = 3346.9333
$ws.Range("I102").Value = 2229
$ws.Range("K102").Value = 2229
$ws.Range("M102").Value = -607
$ws.Range("H109").Value = 50000
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()
$ws.Range("H132").Value = 2759.5715
$ws.Range("I132").Value = 3054.8076
$ws.Range("J132").Value = 1906.6666
$ws.Range("K132").Value = 9164.4228
$ws.Range("L132").Value = 5719.9998
$ws.Range("M132").Value = -6634.4228
$ws.Range("N132").Value = -10779.9998

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1045.625
$ws.Range("J22").Value = 1122.5
$ws.Range("L22").Value = 1122.5
$ws.Range("N22").Value = -1712.5
$ws.Range("H27").Value = 1045.625
$ws.Range("J27").Value = 1122.5
$ws.Range("L27").Value = 1122.5
$ws.Range("N27").Value = -1336.5
$ws.Range("H46").Value = 8157.6875
$ws.Range("I46").Value = 2143.25
$ws.Range("J46").Value = 10162.5
$ws.Range("K46").Value = 2143.25
$ws.Range("L46").Value = 10162.5
$ws.Range("M46").Value = -1955.25
$ws.Range("N46").Value = -10538.5
$ws.Range("H68").Value = 269047.22
$ws.Range("I68").Value = 255000
$ws.Range("J68").Value = 279263.38
$ws.Range("K68").Value = 255000
$ws.Range("L68").Value = 279263.38
$ws.Range("M68").Value = -254251
$ws.Range("N68").Value = -280761.38
$ws.Range("H71").Value = 269047.22
$ws.Range("I71").Value = 255000
$ws.Range("J71").Value = 279263.38
$ws.Range("K71").Value = 1275000
$ws.Range("L71").Value = 1396316.9
$ws.Range("M71").Value = -1271256
$ws.Range("N71").Value = -1403804.9
$ws.Range("H132").Value = 5568.7144
$ws.Range("I132").Value = 5375
$ws.Range("K132").Value = 16125
$ws.Range("M132").Value = -13595

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8773955
$ws.Range("I81").Value = 1862.1666
$ws.Range("J81").Value = 41669304
$ws.Range("K81").Value = 3724.3332
$ws.Range("L81").Value = 83338608
$ws.Range("M81").Value = -2663.3332
$ws.Range("N81").Value = -83340730
$ws.Range("H84").Value = 8773955
$ws.Range("I84").Value = 1862.1666
$ws.Range("J84").Value = 41669304
$ws.Range("K84").Value = 18621.666
$ws.Range("L84").Value = 416693040
$ws.Range("M84").Value = -13317.666
$ws.Range("N84").Value = -416703648
$ws.Range("H109").Value = 27821.846
$ws.Range("I109").Value = 20842
$ws.Range("J109").Value = 29090.908
$ws.Range("K109").Value = 20842
$ws.Range("L109").Value = 29090.908
$ws.Range("M109").Value = -19455
$ws.Range("N109").Value = -31864.908
$ws.Range("H113").Value = 699.59375
$ws.Range("I113").Value = 464.41666
$ws.Range("K113").Value = 1393.24998
$ws.Range("M113").Value = 776.7500199999999
$ws.Range("H132").Value = 3256.84
$ws.Range("I132").Value = 2499
$ws.Range("K132").Value = 7497
$ws.Range("M132").Value = -4967
$ws.Range("H136").Value = 2789.12
$ws.Range("I136").Value = 2178.3171
$ws.Range("K136").Value = 6534.951300000001
$ws.Range("M136").Value = -3984.951300000001

Write-Host "Applied all Ultros_Profits updates"